# Update "想去人数" (want-to-go count) figures on the 展览 (Exhibition)
# and 全部类型 (All Types) sheets, reflecting newer counts scraped at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 1475
$ws1.Range("F3").Value = 3096
$ws1.Range("F4").Value = 43
$ws1.Range("F5").Value = 719
$ws1.Range("F6").Value = 293

# --- Sheet "全部类型" ---
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 1475
$ws4.Range("F3").Value = 3096
$ws4.Range("F4").Value = 43
$ws4.Range("F5").Value = 719
$ws4.Range("F7").Value = 293
